{"js": "// Replace the three-digit-division answer cells with their new values.\n// Each \"old\" string occurs exactly once in the document body, so a plain\n// text search + InsertText(\"Replace\") per pair is sufficient -- no need to\n// depend on table/row/column indices.\nconst replacements = [\n  [\"620\u00f79=68, 8\", \"821\u00f79=91, 2\"],\n  [\"765\u00f75=153, 0\", \"863\u00f74=215, 3\"],\n  [\"755\u00f73=251, 2\", \"828\u00f79=92, 0\"],\n  [\"478\u00f77=68, 2\", \"424\u00f77=60, 4\"],\n  [\"872\u00f74=218, 0\", \"240\u00f72=120, 0\"],\n  [\"357\u00f77=51, 0\", \"925\u00f72=462, 1\"],\n  [\"612\u00f76=102, 0\", \"241\u00f73=80, 1\"],\n  [\"989\u00f76=164, 5\", \"613\u00f72=306, 1\"],\n  [\"860\u00f74=215, 0\", \"272\u00f72=136, 0\"],\n  [\"549\u00f79=61, 0\", \"924\u00f76=154, 0\"],\n  [\"401\u00f73=133, 2\", \"567\u00f75=113, 2\"],\n  [\"250\u00f74=62, 2\", \"833\u00f77=119, 0\"],\n  [\"715\u00f72=357, 1\", \"105\u00f72=52, 1\"],\n  [\"346\u00f75=69, 1\", \"704\u00f79=78, 2\"],\n  [\"906\u00f77=129, 3\", \"978\u00f73=326, 0\"],\n  [\"313\u00f72=156, 1\", \"688\u00f78=86, 0\"],\n  [\"291\u00f75=58, 1\", \"287\u00f76=47, 5\"],\n  [\"566\u00f75=113, 1\", \"259\u00f72=129, 1\"],\n  [\"408\u00f79=45, 3\", \"949\u00f79=105, 4\"],\n  [\"356\u00f77=50, 6\", \"570\u00f74=142, 2\"],\n  [\"530\u00f72=265, 0\", \"990\u00f78=123, 6\"],\n  [\"648\u00f73=216, 0\", \"441\u00f76=73, 3\"],\n  [\"139\u00f75=27, 4\", \"461\u00f78=57, 5\"],\n  [\"670\u00f74=167, 2\", \"587\u00f79=65, 2\"],\n  [\"416\u00f72=208, 0\", \"446\u00f77=63, 5\"]\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text to replace: \"${oldText}\"`);\n  }\n\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the three-digit-division answer cells with the new values.\n# Each old value is a unique string within the document body, so a\n# direct Find/Replace per pair is sufficient and avoids any row/column\n# index assumptions about the table layout.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"620\u00f79=68, 8\", \"821\u00f79=91, 2\"),\n    @(\"765\u00f75=153, 0\", \"863\u00f74=215, 3\"),\n    @(\"755\u00f73=251, 2\", \"828\u00f79=92, 0\"),\n    @(\"478\u00f77=68, 2\", \"424\u00f77=60, 4\"),\n    @(\"872\u00f74=218, 0\", \"240\u00f72=120, 0\"),\n    @(\"357\u00f77=51, 0\", \"925\u00f72=462, 1\"),\n    @(\"612\u00f76=102, 0\", \"241\u00f73=80, 1\"),\n    @(\"989\u00f76=164, 5\", \"613\u00f72=306, 1\"),\n    @(\"860\u00f74=215, 0\", \"272\u00f72=136, 0\"),\n    @(\"549\u00f79=61, 0\", \"924\u00f76=154, 0\"),\n    @(\"401\u00f73=133, 2\", \"567\u00f75=113, 2\"),\n    @(\"250\u00f74=62, 2\", \"833\u00f77=119, 0\"),\n    @(\"715\u00f72=357, 1\", \"105\u00f72=52, 1\"),\n    @(\"346\u00f75=69, 1\", \"704\u00f79=78, 2\"),\n    @(\"906\u00f77=129, 3\", \"978\u00f73=326, 0\"),\n    @(\"313\u00f72=156, 1\", \"688\u00f78=86, 0\"),\n    @(\"291\u00f75=58, 1\", \"287\u00f76=47, 5\"),\n    @(\"566\u00f75=113, 1\", \"259\u00f72=129, 1\"),\n    @(\"408\u00f79=45, 3\", \"949\u00f79=105, 4\"),\n    @(\"356\u00f77=50, 6\", \"570\u00f74=142, 2\"),\n    @(\"530\u00f72=265, 0\", \"990\u00f78=123, 6\"),\n    @(\"648\u00f73=216, 0\", \"441\u00f76=73, 3\"),\n    @(\"139\u00f75=27, 4\", \"461\u00f78=57, 5\"),\n    @(\"670\u00f74=167, 2\", \"587\u00f79=65, 2\"),\n    @(\"416\u00f72=208, 0\", \"446\u00f77=63, 5\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Could not find text to replace: $oldText\"\n    }\n}\n\nWrite-Output \"done\"\n"}
